$d = $word.ActiveDocument
Write-Host "Before: $($d.Paragraphs.Count)"

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Step 1: remove _GoBack bookmark from paragraph 91 (empty paragraph)
$p91 = $d.Paragraphs.Item(91)
$xml91 = "<w:p $ns w:rsidR=`"0040678A`" w:rsidRDefault=`"0040678A`"><w:pPr><w:ind w:left=`"1440`"/><w:contextualSpacing w:val=`"0`"/></w:pPr></w:p>"
$p91.Range.InsertXML($xml91)

# Step 2: locate heading paragraph again (index may be unchanged since same count)
$p92 = $d.Paragraphs.Item(92)
Write-Host "P92 text: [$($p92.Range.Text)]"
$xml92 = "<w:p $ns w:rsidR=`"0040678A`" w:rsidRDefault=`"005841DA`"><w:pPr><w:pStyle w:val=`"Ttulo2`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:bookmarkStart w:id=`"10`" w:name=`"_ebh75tvectrh`" w:colFirst=`"0`" w:colLast=`"0`"/><w:bookmarkStart w:id=`"11`" w:name=`"_Toc525265527`"/><w:bookmarkEnd w:id=`"10`"/><w:r><w:t>Políticas, Directrices y procedimientos</w:t></w:r><w:bookmarkEnd w:id=`"11`"/></w:p>" +
         "<w:p $ns w:rsidR=`"0040678A`" w:rsidRDefault=`"005841DA`"><w:pPr><w:rPr><w:u w:val=`"single`"/></w:rPr></w:pPr><w:r><w:t>Las políticas, directrices y procedimientos a utilizarse son los siguientes:</w:t></w:r><w:bookmarkStart w:id=`"12`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"12`"/></w:p>"
$p92.Range.InsertXML($xml92)

Write-Host "After: $($d.Paragraphs.Count)"
